$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) The "automatic date" placeholder on the slide master and on every
#    slide layout shows a cached value ("23.02.2017"). Bump it forward
#    to "24.04.2017" everywhere it appears.
# ---------------------------------------------------------------------
function UpdateDatePlaceholder($shapes) {
    for ($shpIdx = 1; $shpIdx -le $shapes.Count; $shpIdx++) {
        $shp = $shapes.Item($shpIdx)
        if ($shp.HasTextFrame) {
            $shpRange = $shp.TextFrame.TextRange
            if ($shpRange.Text -eq "23.02.2017") {
                $shpRange.Text = "24.04.2017"
            }
        }
    }
}

UpdateDatePlaceholder $p.SlideMaster.Shapes

for ($layoutIdx = 1; $layoutIdx -le $p.SlideMaster.CustomLayouts.Count; $layoutIdx++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($layoutIdx)
    UpdateDatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 3's title reads "Children" -- extend it with " and more".
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$titleShape = $slide3.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.InsertAfter(" and more")
